# chore(global): update list of issues
# Append two new countries (Honduras, Jamaica) to the "list of issues"
# worksheet, both sharing the existing "Issues with OSM cables" solution
# already used by several other rows (e.g. Armenia, row 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 38: Honduras
$ws.Range("A38").Value = "Honduras"
$ws.Range("B38").Value = "Issues with OSM cables"

# New row 39: Jamaica
$ws.Range("B39").Value = "Issues with OSM cables"
$ws.Range("A39").Value = "Jamaica"

# Match the formatting already used throughout column A/B (wrap text,
# ht=16 row height), as seen on the preceding rows.
$ws.Range("A38:B39").WrapText = $true
$ws.Rows.Item(38).RowHeight = 16
$ws.Rows.Item(39).RowHeight = 16

# Update the saved view/selection state to reflect scrolling down to see
# the newly-added rows (matches the author's saved cursor position).
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("B40").Select()
